$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 168. This shifts the
# existing rows 168-282 down to 169-283 (matching the diff, where every row
# from 168 onward ends up holding the values previously one row above it,
# and a new record appears at the (new) row 168).
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new record's data.
$ws.Cells.Item(168, 1).Value = 5
$ws.Cells.Item(168, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(168, 3).Value = "Maule"
$ws.Cells.Item(168, 4).Value = 44767
$ws.Cells.Item(168, 5).Value = 7
$ws.Cells.Item(168, 6).Value = 100112009
$ws.Cells.Item(168, 7).Value = "Acelga"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 500
$ws.Cells.Item(168, 11).Value = 4000
$ws.Cells.Item(168, 12).Value = 4000
$ws.Cells.Item(168, 13).Value = 4000
$ws.Cells.Item(168, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 1000
$ws.Cells.Item(168, 17).Value = 4
$ws.Cells.Item(168, 18).Value = "Hortaliza"
